# Workflow sheet setup - John Hetrick's commit:
# "Setup and Initialize Workflow creation and processing."
#
# This script:
#  1. Simplifies the MetaData "Seq Key" column (G) so it no longer embeds the
#     sequence number (since the seq number doesn't matter for the key), and
#     documents that with a note in H3.
#  2. Updates the Map sheet key to include /SEQ-0 on the DT-0 key.
#  3. Adds a brand-new "Workflow" sheet listing WFG-<#>/WFS-<#> keys.
#  4. Re-points a couple of view/selection/column-width cosmetics to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) DataSets sheet - select the header rows (cosmetic selection state)
# ---------------------------------------------------------------------
$dataSets = $wb.Worksheets.Item("DataSets")
$dataSets.Range("A1:XFD2").Select()

# ---------------------------------------------------------------------
# 2) MetaData sheet - collapse the Seq Key values, add explanatory note
# ---------------------------------------------------------------------
$metaData = $wb.Worksheets.Item("MetaData")
$metaData.Activate()

$metaData.Columns("H").ColumnWidth = 39.5859375

$metaData.Range("G2").Value = "INST-1/EXP-1/GRP-0/DT-0"
$metaData.Range("G3").Value = "INST-1/EXP-1/GRP-0/DT-0"
$metaData.Range("G4").Value = "INST-1/EXP-1/GRP-0/DT-0"

$metaData.Range("G5").Value = "INST-1/EXP-1/GRP-1/DT-2"
$metaData.Range("G6").Value = "INST-1/EXP-1/GRP-1/DT-2"
$metaData.Range("G7").Value = "INST-1/EXP-1/GRP-1/DT-2"

$metaData.Range("G8").Value = "INST-2/EXP-2/GRP-1/DT-0"

$noteCell = $metaData.Range("H3")
$noteCell.Value = "Note the Seq. number is not important for the key.  All sequences use the same workflow description"
$noteCell.Borders.Item(7).LineStyle = 1
$noteCell.Borders.Item(7).Weight = 2
$noteCell.Borders.Item(10).LineStyle = 1
$noteCell.Borders.Item(10).Weight = 2

$metaData.Range("G2").Select()

# ---------------------------------------------------------------------
# 3) Map sheet - the Group/DataType key now also carries the Seq
# ---------------------------------------------------------------------
$map = $wb.Worksheets.Item("Map")
$map.Activate()
$map.Columns("A").ColumnWidth = 34.29296875
$map.Range("A3").Value = "INST-1/EXP-1/GRP-0/DT-0/SEQ-0"
$map.Range("A3").Select()

# ---------------------------------------------------------------------
# 4) Brand new "Workflow" sheet, appended after "Map"
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$workflow = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$workflow.Name = "Workflow"

$workflow.Columns("A").ColumnWidth = 28.5859375
$workflow.Columns("B").ColumnWidth = 22.46875

$workflow.Range("A1").Value = "ID"
$workflow.Range("A1").Style = $metaData.Range("A1").Style

$workflow.Range("A2").Value = "INST-1/EXP-1/GRP-0/WFG-<#>"
$workflow.Range("A2").Style = $metaData.Range("A2").Style

$workflow.Range("A3").Value = "INST-1/EXP-1/GRP-0/WFS-<#>"
$workflow.Range("A3").Style = $metaData.Range("A2").Style

$workflow.Range("B1:B3").Style = $workflow.Range("A1").Style

$workflow.Range("C5").Select()
